$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 132356.73
$ws.Range("I15").Value = 132356.73
$ws.Range("K15").Value = 397070.1900000001
$ws.Range("M15").Value = -396901.1900000001
$ws.Range("H17").Value = 783137.4
$ws.Range("J17").Value = 916394.8
$ws.Range("L17").Value = 2749184.4
$ws.Range("N17").Value = -2749520.4
$ws.Range("H33").Value = 139.77777
$ws.Range("I33").Value = 136.14285
$ws.Range("K33").Value = 136.14285
$ws.Range("M33").Value = 92.85714999999999
$ws.Range("H69").Value = 4631829.5
$ws.Range("I69").Value = 2045.4546
$ws.Range("J69").Value = 8549339
$ws.Range("K69").Value = 6136.3638
$ws.Range("L69").Value = 25648017
$ws.Range("M69").Value = -5262.3638
$ws.Range("N69").Value = -25649765
$ws.Range("H72").Value = 4631829.5
$ws.Range("I72").Value = 2045.4546
$ws.Range("J72").Value = 8549339
$ws.Range("K72").Value = 18409.0914
$ws.Range("L72").Value = 76944051
$ws.Range("M72").Value = -14041.0914
$ws.Range("N72").Value = -76952787
$ws.Range("H105").Value = 500335.5
$ws.Range("J105").Value = 500335.5
$ws.Range("L105").Value = 500335.5
$ws.Range("N105").Value = -507323.5
$ws.Range("H132").Value = 247947.64
$ws.Range("I132").Value = 259510.77
$ws.Range("J132").Value = 66792
$ws.Range("K132").Value = 778532.3099999999
$ws.Range("L132").Value = 200376
$ws.Range("M132").Value = -776002.3099999999
$ws.Range("N132").Value = -205436
$ws.Range("H133").Value = 12000
$ws.Range("J133").Value = 12000
$ws.Range("L133").Value = 12000
$ws.Range("N133").Value = -22120
$ws.Range("H135").Value = 1287.8485
$ws.Range("I135").Value = 1128.52
$ws.Range("J135").Value = 1785.75
$ws.Range("K135").Value = 10156.68
$ws.Range("L135").Value = 16071.75
$ws.Range("M135").Value = -7621.68
$ws.Range("N135").Value = -21141.75
$ws.Range("H138").Value = 4597260.5
$ws.Range("I138").Value = 837758.9
$ws.Range("J138").Value = 11113730
$ws.Range("K138").Value = 2513276.7
$ws.Range("L138").Value = 33341190
$ws.Range("M138").Value = -2508136.7
$ws.Range("N138").Value = -33351470

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1519
$ws.Range("I61").Value = 1519
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1519
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1307
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 4342.1816
$ws.Range("I74").Value = 1414.258
$ws.Range("J74").Value = 11324.154
$ws.Range("K74").Value = 1414.258
$ws.Range("L74").Value = 11324.154
$ws.Range("M74").Value = -540.258
$ws.Range("N74").Value = -13072.154
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 4342.1816
$ws.Range("I77").Value = 1414.258
$ws.Range("J77").Value = 11324.154
$ws.Range("K77").Value = 7071.29
$ws.Range("L77").Value = 56620.77
$ws.Range("M77").Value = -2703.29
$ws.Range("N77").Value = -65356.77
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H132").Value = 2172.1042
$ws.Range("I132").Value = 1836.0952
$ws.Range("J132").Value = 4524.1665
$ws.Range("K132").Value = 5508.2856
$ws.Range("L132").Value = 13572.4995
$ws.Range("M132").Value = -2978.2856
$ws.Range("N132").Value = -18632.4995
$ws.Range("H136").Value = 1519
$ws.Range("I136").Value = 1519
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4557
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2007
$ws.Range("N136").ClearContents()
$ws.Range("H139").Value = 48163.75
$ws.Range("I139").Value = 30470
$ws.Range("J139").Value = 65857.5
$ws.Range("K139").Value = 30470
$ws.Range("L139").Value = 65857.5
$ws.Range("M139").Value = -25330
$ws.Range("N139").Value = -76137.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 28574744
$ws.Range("I134").Value = 38464144
$ws.Range("K134").Value = 115392432
$ws.Range("M134").Value = -115389897

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1407.9354
$ws.Range("I58").Value = 742.1818
$ws.Range("J58").Value = 3035.3333
$ws.Range("K58").Value = 742.1818
$ws.Range("L58").Value = 3035.3333
$ws.Range("M58").Value = -539.1818
$ws.Range("N58").Value = -3441.3333
$ws.Range("H105").Value = 848.36365
$ws.Range("I105").Value = 849.7368
$ws.Range("J105").Value = 839.6667
$ws.Range("K105").Value = 849.7368
$ws.Range("L105").Value = 839.6667
$ws.Range("M105").Value = 897.2632
$ws.Range("N105").Value = -4333.6667
$ws.Range("H132").Value = 1673.5957
$ws.Range("I132").Value = 1532.5366
$ws.Range("J132").Value = 2637.5
$ws.Range("K132").Value = 4597.6098
$ws.Range("L132").Value = 7912.5
$ws.Range("M132").Value = -2067.6098
$ws.Range("N132").Value = -12972.5
$ws.Range("H134").Value = 2133.85
$ws.Range("I134").Value = 1392.262
$ws.Range("J134").Value = 3864.2222
$ws.Range("K134").Value = 4176.786
$ws.Range("L134").Value = 11592.6666
$ws.Range("M134").Value = -1641.786
$ws.Range("N134").Value = -16662.6666
$ws.Range("H136").Value = 1407.9354
$ws.Range("I136").Value = 742.1818
$ws.Range("J136").Value = 3035.3333
$ws.Range("K136").Value = 2226.5454
$ws.Range("L136").Value = 9105.999899999999
$ws.Range("M136").Value = 323.4546
$ws.Range("N136").Value = -14205.9999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 26316466
$ws.Range("I113").Value = 713
$ws.Range("J113").Value = 35714948
$ws.Range("K113").Value = 2139
$ws.Range("L113").Value = 107144844
$ws.Range("M113").Value = 31
$ws.Range("N113").Value = -107149184
$ws.Range("H131").Value = 1894.238
$ws.Range("I131").Value = 482.5
$ws.Range("J131").Value = 2226.4119
$ws.Range("K131").Value = 1447.5
$ws.Range("L131").Value = 6679.2357
$ws.Range("M131").Value = 3592.5
$ws.Range("N131").Value = -16759.2357

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2974.8293
$ws.Range("I132").Value = 2769.027
$ws.Range("J132").Value = 4878.5
$ws.Range("K132").Value = 8307.081
$ws.Range("L132").Value = 14635.5
$ws.Range("M132").Value = -5777.081
$ws.Range("N132").Value = -19695.5
$ws.Range("H137").Value = 54766.668
$ws.Range("J137").Value = 54766.668
$ws.Range("L137").Value = 54766.668
$ws.Range("N137").Value = -64966.668

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 32681
$ws.Range("J81").Value = 32681
$ws.Range("L81").Value = 32681
$ws.Range("N81").Value = -34677
$ws.Range("H82").Value = 1600
$ws.Range("I82").Value = 1600
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1600
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1239
$ws.Range("N82").ClearContents()
$ws.Range("H84").Value = 32681
$ws.Range("J84").Value = 32681
$ws.Range("L84").Value = 98043
$ws.Range("N84").Value = -108027
$ws.Range("H85").Value = 1600
$ws.Range("I85").Value = 1600
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1600
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -352
$ws.Range("N85").ClearContents()
$ws.Range("H132").Value = 4156.054
$ws.Range("I132").Value = 4119
$ws.Range("J132").Value = 4263.316
$ws.Range("K132").Value = 12357
$ws.Range("L132").Value = 12789.948
$ws.Range("M132").Value = -9827
$ws.Range("N132").Value = -17849.948
$ws.Range("H136").Value = 4574.227
$ws.Range("I136").Value = 3017.889
$ws.Range("J136").Value = 11577.75
$ws.Range("K136").Value = 9053.667000000001
$ws.Range("L136").Value = 34733.25
$ws.Range("M136").Value = -6503.667000000001
$ws.Range("N136").Value = -39833.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1252452.5
$ws.Range("I81").Value = 1540472.2
$ws.Range("J81").Value = 4366.6665
$ws.Range("K81").Value = 3080944.4
$ws.Range("L81").Value = 8733.333000000001
$ws.Range("M81").Value = -3079883.4
$ws.Range("N81").Value = -10855.333
$ws.Range("H84").Value = 1252452.5
$ws.Range("I84").Value = 1540472.2
$ws.Range("J84").Value = 4366.6665
$ws.Range("K84").Value = 15404722
$ws.Range("L84").Value = 43666.665
$ws.Range("M84").Value = -15399418
$ws.Range("N84").Value = -54274.665
$ws.Range("H132").Value = 7938187
$ws.Range("I132").Value = 9260776
$ws.Range("J132").Value = 2649.889
$ws.Range("K132").Value = 27782328
$ws.Range("L132").Value = 7949.667
$ws.Range("M132").Value = -27779798
$ws.Range("N132").Value = -13009.667
$ws.Range("H136").Value = 17817.283
$ws.Range("I136").Value = 18533.719
$ws.Range("K136").Value = 55601.15700000001
$ws.Range("M136").Value = -53051.15700000001

Write-Host "Applied all updates"